# Update the "取得日時" (timestamp) column (A) for all data rows (2-16)
# on the "ランサーズ" worksheet to the new append timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-09-26 12:50:21"

for ($row = 2; $row -le 16; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
